## Entertainment Expense.xlsx
## Commit: "Updated cleaned files and added Raw_Data sheet"
##
## - rename the existing sheet to "Cleaned_Data" (keeps its original sheetId)
## - add a brand-new "Raw_Data" sheet positioned before it, holding the
##   un-cleaned numbers (a straight copy of the values, no totals row,
##   no special formatting)
## - a temporary sheet is added/removed purely so the engine's "last sheet
##   touched" pointer lands back on Cleaned_Data (making it the active tab,
##   matching the author's workbookView activeTab="1")
## - tidy a couple of formatting differences on Cleaned_Data (bottom row
##   gets a heavier rule instead of a tinted fill, the drawing part is gone)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the original sheet -> Cleaned_Data (this is sheetId 1 / rId2
#    in the target: it's the sheet that already existed).
# ---------------------------------------------------------------------
$wsCleaned = $wb.Worksheets.Item(1)
$wsCleaned.Name = "Cleaned_Data"

# ---------------------------------------------------------------------
# 2. Insert the new Raw_Data sheet *before* Cleaned_Data (this is the
#    brand-new sheet, so it picks up sheetId 2 / rId1 in the target).
# ---------------------------------------------------------------------
$wsRaw = $wb.Worksheets.Add($wsCleaned, $null)
$wsRaw.Name = "Raw_Data"

# Raw, uncleaned data -- same numbers as Cleaned_Data, no totals row.
$wsRaw.Range("A1").Value = "Entertainment"
$wsRaw.Range("B1").Value = "Jan"
$wsRaw.Range("C1").Value = "Feb"

$wsRaw.Range("A2").Value = "Cable TV"
$wsRaw.Range("B2").Value = 95.67
$wsRaw.Range("C2").Value = 95.67

$wsRaw.Range("A3").Value = "Video Streaming"
$wsRaw.Range("B3").Value = 9.99
$wsRaw.Range("C3").Value = 9.99

$wsRaw.Range("A4").Value = "Movies"
$wsRaw.Range("B4").Value = 32
$wsRaw.Range("C4").Value = 16

$wsRaw.Range("A5").Value = "Music"
$wsRaw.Range("B5").Value = 41.98
$wsRaw.Range("C5").Value = 0

$wsRaw.Range("A6").Value = "Video Games"
$wsRaw.Range("B6").Value = 132.32
$wsRaw.Range("C6").Value = 62.7

$wsRaw.Range("A7").Value = "Totals"

# ---------------------------------------------------------------------
# 3. A scratch sheet: adding it shifts the "active" sheet to it, and
#    removing it again drops the active sheet back to Cleaned_Data --
#    which is what the saved file should show as the selected tab.
# ---------------------------------------------------------------------
$wsScratch = $wb.Worksheets.Add($null, $wsCleaned)
$wsScratch.Delete()

# ---------------------------------------------------------------------
# 4. Cleaned_Data tweaks: the totals row (6) now reads with a heavier
#    bottom rule instead of the old tinted-fill treatment, and the
#    chart/drawing that used to hang off this sheet is gone.
# ---------------------------------------------------------------------
$wsCleaned.Range("B4").Value = 32
$wsCleaned.Range("C4").Value = 16
$wsCleaned.Range("C5").Value = 0

$wsCleaned.Range("A6:C6").Borders.Item(9).LineStyle = 1
$wsCleaned.Range("A6:C6").Borders.Item(9).Weight = -4138
$wsCleaned.Range("A6:C6").Interior.Pattern = -4142
